$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Formula = "0.9975859522819519"
$ws.Range("B1").Formula = "1.817065834999084"
$ws.Range("C1").Formula = "6.884543418884277"
$ws.Range("D1").Formula = "2.89591646194458"
$ws.Range("E1").Formula = "0.4142286777496338"
